# Auto-generated script applying the cryptos.xlsx price/volume/coin-listing update
# Commit: Updated symbol list on Sat Jan 28 16:40:06 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Addr, $Text)
    $cell = $Sheet.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

# Row 2
Set-TextCell $ws "D2" "306.98"
Set-TextCell $ws "E2" "0.55%"

# Row 3
Set-TextCell $ws "D3" "38.97"
Set-TextCell $ws "E3" "8.34%"

# Row 4
Set-TextCell $ws "D4" "5.094"
Set-TextCell $ws "E4" "0.90%"

# Row 5
Set-TextCell $ws "D5" "0.08059"
Set-TextCell $ws "E5" "0.18%"

# Row 6
Set-TextCell $ws "D6" "1.919"
Set-TextCell $ws "E6" "2.69%"

# Row 7
Set-TextCell $ws "B7" "KuCoinToken"
Set-TextCell $ws "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell $ws "D7" "7.940"
Set-TextCell $ws "E7" "1.98%"

# Row 8
Set-TextCell $ws "B8" "MXToken"
Set-TextCell $ws "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D8" "0.9305"
Set-TextCell $ws "E8" "0.44%"

# Row 9
Set-TextCell $ws "B9" "LiechtensteinCryptoassetsExchange"
Set-TextCell $ws "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell $ws "D9" "0.1454"
Set-TextCell $ws "E9" "6.67%"

# Row 10
Set-TextCell $ws "B10" "WazirX"
Set-TextCell $ws "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws "D10" "0.1935"
Set-TextCell $ws "E10" "2.03%"

# Row 11
Set-TextCell $ws "B11" "MandalaExchangeToken"
Set-TextCell $ws "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell $ws "D11" "0.09011"
Set-TextCell $ws "E11" "-0.47%"

# Row 12
Set-TextCell $ws "B12" "BitrueCoin"
Set-TextCell $ws "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell $ws "D12" "0.03508"
Set-TextCell $ws "E12" "2.24%"

# Row 13
Set-TextCell $ws "B13" "BitMartToken"
Set-TextCell $ws "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell $ws "D13" "0.09786"
Set-TextCell $ws "E13" "-1.06%"

# Row 14
Set-TextCell $ws "B14" "BitForexToken"
Set-TextCell $ws "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws "D14" "0.001397"
Set-TextCell $ws "E14" "-0.48%"

# Row 15
Set-TextCell $ws "B15" "TigerCash"
Set-TextCell $ws "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D15" "0.006114"
Set-TextCell $ws "E15" "0.41%"

# Row 16
Set-TextCell $ws "B16" "LEO"
Set-TextCell $ws "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D16" "3.755"
Set-TextCell $ws "E16" "-2.15%"

# Row 17
Set-TextCell $ws "B17" "GateToken"
Set-TextCell $ws "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws "D17" "4.190"
Set-TextCell $ws "E17" "1.71%"

# Row 18
Set-TextCell $ws "E18" "2.19%"

# Row 19
Set-TextCell $ws "D19" "0.3444"
Set-TextCell $ws "E19" "0.89%"

# Row 20
Set-TextCell $ws "D20" "0.1319"
Set-TextCell $ws "E20" "1.49%"

# Row 21
Set-TextCell $ws "D21" "4.781"
Set-TextCell $ws "E21" "-1.00%"

# Row 22
Set-TextCell $ws "D22" "0.2504"
Set-TextCell $ws "E22" "5.03%"

# Row 23
Set-TextCell $ws "D23" "0.04379"
Set-TextCell $ws "E23" "0.47%"

# Row 24
Set-TextCell $ws "D24" "0.001235"
Set-TextCell $ws "E24" "0.46%"

# Row 25
Set-TextCell $ws "D25" "0.004279"
Set-TextCell $ws "E25" "-0.17%"

# Row 26
Set-TextCell $ws "D26" "0.0001300"
Set-TextCell $ws "E26" "-0.01%"

# Row 39
Set-TextCell $ws "D39" "0.02070"
Set-TextCell $ws "E39" "3.48%"

# Row 40
Set-TextCell $ws "D40" "0.05050"
Set-TextCell $ws "E40" "-1.07%"

# Row 41
Set-TextCell $ws "D41" "0.007437"
Set-TextCell $ws "E41" "-1.03%"

# Row 42
Set-TextCell $ws "D42" "0.01007"
Set-TextCell $ws "E42" "-0.04%"

# Row 43
Set-TextCell $ws "D43" "0.1351"
Set-TextCell $ws "E43" "-0.43%"

# Row 44
Set-TextCell $ws "D44" "0.002141"
Set-TextCell $ws "E44" "-0.94%"

# Row 45
Set-TextCell $ws "D45" "0.009077"
Set-TextCell $ws "E45" "-5.63%"

# Row 46
Set-TextCell $ws "D46" "0.00006188"
Set-TextCell $ws "E46" "-0.53%"

# Row 47
Set-TextCell $ws "E47" "0.14%"

# Row 48
Set-TextCell $ws "D48" "0.002798"

# Row 49
Set-TextCell $ws "D49" "0.001599"
Set-TextCell $ws "E49" "28.09%"

# Row 50
Set-TextCell $ws "D50" "0.00002101"
Set-TextCell $ws "E50" "0.14%"

# Row 51
Set-TextCell $ws "D51" "0.0002001"
Set-TextCell $ws "E51" "0.14%"
